$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (rows 2 and 3), pushing
# everything else down by two rows.
$ws.Rows("2:3").Insert()

# Row 2: newest survey (Survey 44), baseline only.
# Row 3: Survey 43 follow-up, with modules/topics note.
$ws.Range("A2").Value = "Survey 44"
$ws.Range("A3").Value = "Survey 43"
$ws.Range("B3").Value = "Feb 2 - Feb 4"
$ws.Range("B2").Value = "Feb 9 - Feb 11"
$ws.Range("C2").Value = "x"
$ws.Range("D3").Value = "x"
$ws.Range("E3").Value = "Social Support`nHealthcare`nVaccines`nRISER`nUnemployment`nStimulus"

# The row insert carried an empty styled cell into E2 (from the header
# row's formatting); clear it so row 2 has no stray E cell.
$ws.Range("E2").Clear()

# Match the autofit height Excel would compute for the wrapped 6-line note.
$ws.Range("E3").EntireRow.RowHeight = 102

# Update selection/view to match the saved state (E5 selected, no pinned
# scroll position).
$ws.Range("E5").Select()
